$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply the same number formats used by the existing data rows so the new
# rows reuse the existing style indices (date / 0.00 / integer formats)
# instead of creating new ones in styles.xml.
$dateFmt = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'
$decimalFmt = '0.00_);[Red]\(0.00\)'
$intFmt = '0_);[Red]\(0\)'

$ws.Range("A54:A55").NumberFormat = $dateFmt
$ws.Range("C54:E55").NumberFormat = $decimalFmt
$ws.Range("F54:F55").NumberFormat = $intFmt

# Row 54: 2025-09-27, 四方坪站
$ws.Range("A54").Value = 45927
$ws.Range("B54").Value = "四方坪站"
$ws.Range("C54").Value = 9689.46
$ws.Range("D54").Value = 8028.25
$ws.Range("E54").Value = 3285.07
$ws.Range("F54").Value = 408

# Row 55: 2025-09-27, 高岭站
$ws.Range("A55").Value = 45927
$ws.Range("B55").Value = "高岭站"
$ws.Range("C55").Value = 4148.7
$ws.Range("D55").Value = 3342.83
$ws.Range("E55").Value = 995.63
$ws.Range("F55").Value = 156

$ws.Range("H44").Activate()
